# Add a new "Jobs" worksheet (copied from "Login") with its own set of
# test-suite rows, matching the target diff:
#  - workbook.xml gets a second <sheet> entry and activeTab="1"
#  - sharedStrings.xml gains 4 new strings (SafeWay_Jobs1..4)
#  - sheet1 (Login) loses tabSelected and keeps selection A2
#  - sheet2 (Jobs) is new, tabSelected, selection A18, with updated
#    values/styles on the 4 "section header" rows (3, 8, 13, 17)

$wb = $excel.ActiveWorkbook

$login = $wb.Worksheets.Item("Login")

# Duplicate the Login sheet right after itself; Excel automatically makes
# the new copy the active/selected tab (matching tabSelected moving over).
$login.Copy($null, $login) | Out-Null

$jobs = $wb.Worksheets.Item(2)
$jobs.Name = "Jobs"

# Update the four section-header cells with the new "Jobs" titles.
$jobs.Range("A3").Value = "SafeWay_Jobs1"
$jobs.Range("A8").Value = "SafeWay_Jobs2"
$jobs.Range("A13").Value = "SafeWay_Jobs3"
$jobs.Range("A17").Value = "SafeWay_Jobs4"

# Those same rows change from the "orange" section style to the plain
# bordered style - copy that formatting over from a row that already
# uses it (row 4 uses cell style 6: border only, no fill).
$jobs.Range("A4").Copy() | Out-Null
$jobs.Range("A3").PasteSpecial(-4122) | Out-Null
$jobs.Range("A8").PasteSpecial(-4122) | Out-Null
$jobs.Range("A13").PasteSpecial(-4122) | Out-Null
$jobs.Range("A17").PasteSpecial(-4122) | Out-Null

# Final selection/active sheet state for Jobs.
$jobs.Activate() | Out-Null
$jobs.Range("A18").Select() | Out-Null
